$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (tab name / workbook.xml sheet name)
$ws.Name = "Through 2022-10-31"

# Update the "October" row label text
$ws.Range("A11").Value = "October (through 10-31)"

# Update October row (row 11) values
$ws.Range("B11").Value = 32
$ws.Range("C11").Value = 57
$ws.Range("D11").Value = 83
$ws.Range("E11").Value = 67
$ws.Range("F11").Value = 60
$ws.Range("G11").Value = 156
$ws.Range("H11").Value = 194
$ws.Range("I11").Value = 125

# Update Total row (row 12) values
$ws.Range("B12").Value = 258
$ws.Range("C12").Value = 486
$ws.Range("D12").Value = 710
$ws.Range("E12").Value = 615
$ws.Range("F12").Value = 482
$ws.Range("G12").Value = 1057
$ws.Range("H12").Value = 1441
$ws.Range("I12").Value = 1401
